# Update "docs/epexspot_prices.xlsx" with the newest day of data for each
# sheet: an extra "21-aug" day column on "Prix Spot" (24 hourly values),
# and an extra "2025-08-19" row on "Gaz" and "CO2".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append column BQ (69) = "21-aug" with its 24 hourly
# values, copying the header cell's format from BP1 (the previous day).
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("BP1").Copy()
$wsPrix.Range("BQ1").PasteSpecial(-4122)  # xlPasteFormats
$wsPrix.Cells.Item(1, 69).Value = "21-aug"

$bqValues = @{
    2  = 26.33
    3  = 19.37
    4  = 44.65
    5  = 37.99
    6  = 32.93
    7  = 16.81
    8  = 39.93
    9  = 40.2
    10 = 61.06
    11 = 50.2
    12 = 25.67
    13 = 10
    14 = 7.89
    15 = 5.11
    16 = 5.11
    17 = 3.52
    18 = 5.79
    19 = 8.99
    20 = 25
    21 = 50
    22 = 77.84
    23 = 83.26000000000001
    24 = 55.77
    25 = 52.05
}

foreach ($row in $bqValues.Keys) {
    $wsPrix.Cells.Item($row, 69).Value = $bqValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 66 = 2025-08-19 / 29.8
#
# Note: the existing date column stores each date as a plain text literal
# (not an Excel date serial). Assigning the literal string straight to
# `.Value` gets auto-recognised as a date by Excel's input parser, which
# would both change the stored type and mint a brand-new date-formatted
# style. Routing it through a formula that yields a string, then collapsing
# that formula down to its static value with Paste Special > Values, keeps
# the literal text (and the default, unstyled cell) exactly like the rest
# of the column.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Cells.Item(66, 1).Formula = '="2025-08-19"'
$wsGaz.Range("A66").Copy()
$wsGaz.Range("A66").PasteSpecial(-4163)  # xlPasteValues
$wsGaz.Cells.Item(66, 2).Value = 29.8

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 66 = 2025-08-19 / 71.3 (same text-literal trick)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Cells.Item(66, 1).Formula = '="2025-08-19"'
$wsCo2.Range("A66").Copy()
$wsCo2.Range("A66").PasteSpecial(-4163)  # xlPasteValues
$wsCo2.Cells.Item(66, 2).Value = 71.3

$excel.CutCopyMode = $false
